$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Artist names (column A) were re-shuffled / some swapped for new acts,
#    and view counts (column B) were bumped up ("Ajeitei os dados de lives").
# ---------------------------------------------------------------------------

$ws.Range("A6").Value  = "Raça Negra"
$ws.Range("A7").Value  = "Alok"
$ws.Range("A8").Value  = "Zé Neto e Cristiano"
$ws.Range("A9").Value  = "Jorge e Mateus"
$ws.Range("A12").Value = "Thiaguinho"
$ws.Range("A13").Value = "Maiara e Maraisa"
$ws.Range("A14").Value = " Os Barões da Pisadinha"
$ws.Range("A15").Value = "Marcos e Belutti"
$ws.Range("A17").Value = "Wesley Safadão"
$ws.Range("A18").Value = "Ferrugem"
$ws.Range("A19").Value = "Xand Avião"
$ws.Range("A20").Value = "Gusttavo Lima"
$ws.Range("A21").Value = "César Menotti e Fabiano"

$ws.Range("B2").Value  = 55000000
$ws.Range("B7").Value  = 15000000
$ws.Range("B8").Value  = 14000000
$ws.Range("B9").Value  = 13000000
$ws.Range("B11").Value = 10000000
$ws.Range("B12").Value = 10000000
$ws.Range("B13").Value = 9000000
$ws.Range("B14").Value = 8300000
$ws.Range("B15").Value = 8000000
$ws.Range("B16").Value = 7800000
$ws.Range("B17").Value = 6700000
$ws.Range("B18").Value = 6700000
$ws.Range("B19").Value = 6200000
$ws.Range("B20").Value = 6200000
$ws.Range("B21").Value = 6000000

# ---------------------------------------------------------------------------
# 2) Formatting tweaks that came along with the data edit: B12 right aligned,
#    A15/B15 got a pasted-in font (black Arial) with left/right + vertical
#    center alignment.
# ---------------------------------------------------------------------------

$ws.Range("B12").HorizontalAlignment = -4152

$a15 = $ws.Range("A15")
$a15.Font.Color = 0
$a15.Font.Name = "arial"
$a15.HorizontalAlignment = -4131
$a15.VerticalAlignment = -4108

$b15 = $ws.Range("B15")
$b15.Font.Color = 0
$b15.Font.Name = "Arial"
$b15.HorizontalAlignment = -4152
$b15.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 3) View-state tweaks captured in the sheetView (zoom + scroll position).
# ---------------------------------------------------------------------------

$ws.Application.ActiveWindow.Zoom = 140
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A24:B24").Select()
